$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the single-line mailing address into two lines (street, then city/state/zip)
#    and add a trailing blank line, using a wildcard-free Find/Replace with paragraph marks.
$d.Content.Find.Execute("4638 Lockridge Way, Castro Valley CA 94546", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "4638 Lockridge Way^pCastro Valley, CA 94546^p", 2) | Out-Null

# 3. Remove the two blank paragraphs that directly follow the
#    "...Board of Directors" signature line.
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -match "Board of Directors") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $delIndex = $targetIndex + 1

    $p1 = $d.Paragraphs.Item($delIndex)
    $p1.Range.Delete()

    $p2 = $d.Paragraphs.Item($delIndex)
    $p2.Range.Delete()
}
